# Applies the fixture-data refresh described in the commit:
# "Atualizacao de bases das ligas, do dia: 31-03-2024 as 20:29"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 236-239: correct fixture order & refresh odds (only changed cells touched) ---
# Row 236
$ws.Cells.Item(236, 2).Value = 6861095
$ws.Cells.Item(236, 6).Value = 'FC Botosani'
$ws.Cells.Item(236, 7).Value = 'Farul Constanta'
$ws.Cells.Item(236, 8).Value = 0
$ws.Cells.Item(236, 10).Value = 'D'
$ws.Cells.Item(236, 11).Value = 3.75
$ws.Cells.Item(236, 13).Value = 1.909
$ws.Cells.Item(236, 14).Value = 3.1
$ws.Cells.Item(236, 15).Value = 3
$ws.Cells.Item(236, 16).Value = 2.375
$ws.Cells.Item(236, 17).Value = 0.25
$ws.Cells.Item(236, 18).Value = 1.775
$ws.Cells.Item(236, 19).Value = 2.1
$ws.Cells.Item(236, 20).Value = 2
$ws.Cells.Item(236, 21).Value = 1.8
$ws.Cells.Item(236, 22).Value = 2.05
$ws.Cells.Item(236, 23).Value = -1
$ws.Cells.Item(236, 24).Value = 2
$ws.Cells.Item(236, 26).Value = 0.3875
$ws.Cells.Item(236, 29).Value = 1.05

# Row 237
$ws.Cells.Item(237, 2).Value = 6865915
$ws.Cells.Item(237, 6).Value = 'FC Voluntari'
$ws.Cells.Item(237, 7).Value = 'Universitatea Cluj'
$ws.Cells.Item(237, 8).Value = 0
$ws.Cells.Item(237, 9).Value = 0
$ws.Cells.Item(237, 10).Value = 'D'
$ws.Cells.Item(237, 11).Value = 3.5
$ws.Cells.Item(237, 12).Value = 3.25
$ws.Cells.Item(237, 13).Value = 2.05
$ws.Cells.Item(237, 14).Value = 3.4
$ws.Cells.Item(237, 15).Value = 3.1
$ws.Cells.Item(237, 16).Value = 2.15
$ws.Cells.Item(237, 18).Value = 1.975
$ws.Cells.Item(237, 19).Value = 1.875
$ws.Cells.Item(237, 21).Value = 2.05
$ws.Cells.Item(237, 22).Value = 1.75
$ws.Cells.Item(237, 24).Value = 2.1
$ws.Cells.Item(237, 25).Value = -1
$ws.Cells.Item(237, 26).Value = 0.4875
$ws.Cells.Item(237, 27).Value = -0.5
$ws.Cells.Item(237, 28).Value = -1
$ws.Cells.Item(237, 29).Value = 0.75

# Row 238
$ws.Cells.Item(238, 2).Value = 6870268
$ws.Cells.Item(238, 6).Value = 'Petrolul Ploiesti'
$ws.Cells.Item(238, 7).Value = 'ACS Sepsi'
$ws.Cells.Item(238, 8).Value = 1
$ws.Cells.Item(238, 9).Value = 2
$ws.Cells.Item(238, 10).Value = 'A'
$ws.Cells.Item(238, 11).Value = 2.8
$ws.Cells.Item(238, 12).Value = 3
$ws.Cells.Item(238, 13).Value = 2.55
$ws.Cells.Item(238, 14).Value = 3
$ws.Cells.Item(238, 15).Value = 3.2
$ws.Cells.Item(238, 16).Value = 2.3
$ws.Cells.Item(238, 18).Value = 1.85
$ws.Cells.Item(238, 19).Value = 2
$ws.Cells.Item(238, 20).Value = 2.25
$ws.Cells.Item(238, 21).Value = 1.875
$ws.Cells.Item(238, 22).Value = 1.975
$ws.Cells.Item(238, 24).Value = -1
$ws.Cells.Item(238, 25).Value = 1.3
$ws.Cells.Item(238, 26).Value = -1
$ws.Cells.Item(238, 27).Value = 1
$ws.Cells.Item(238, 28).Value = 0.875
$ws.Cells.Item(238, 29).Value = -1

# Row 239
$ws.Cells.Item(239, 2).Value = 6836277
$ws.Cells.Item(239, 6).Value = 'CFR Cluj'
$ws.Cells.Item(239, 7).Value = 'AFC Hermannstadt'
$ws.Cells.Item(239, 8).Value = 1
$ws.Cells.Item(239, 10).Value = 'H'
$ws.Cells.Item(239, 11).Value = 1.7
$ws.Cells.Item(239, 12).Value = 3.4
$ws.Cells.Item(239, 13).Value = 5
$ws.Cells.Item(239, 14).Value = 1.65
$ws.Cells.Item(239, 15).Value = 3.5
$ws.Cells.Item(239, 16).Value = 5.25
$ws.Cells.Item(239, 17).Value = -0.75
$ws.Cells.Item(239, 18).Value = 1.85
$ws.Cells.Item(239, 19).Value = 2
$ws.Cells.Item(239, 21).Value = 1.875
$ws.Cells.Item(239, 22).Value = 1.975
$ws.Cells.Item(239, 23).Value = 0.6499999999999999
$ws.Cells.Item(239, 24).Value = -1
$ws.Cells.Item(239, 26).Value = 0.425
$ws.Cells.Item(239, 29).Value = 0.9750000000000001
# --- Rows 251-253: fully refresh with new fixtures (old ids superseded) ---
# Clear old contents first (old rows had fewer populated columns; ensure stale H/I/J etc. do not linger)
$ws.Range("A251:AC253").ClearContents()

# --- Rows 254-256 are brand new: clone number-format/style from row 253 (A col) and (E col) ---
$ws.Cells.Item(253, 1).Copy() | Out-Null
$ws.Cells.Item(254, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(253, 5).Copy() | Out-Null
$ws.Cells.Item(254, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(253, 1).Copy() | Out-Null
$ws.Cells.Item(255, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(253, 5).Copy() | Out-Null
$ws.Cells.Item(255, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(253, 1).Copy() | Out-Null
$ws.Cells.Item(256, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(253, 5).Copy() | Out-Null
$ws.Cells.Item(256, 5).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Write full data for rows 251-256 ---
# Row 251
$ws.Cells.Item(251, 1).Value = 249
$ws.Cells.Item(251, 2).Value = 7951781
$ws.Cells.Item(251, 3).Value = 'Romania Liga I'
$ws.Cells.Item(251, 4).Value = 'Romania Liga I'
$ws.Cells.Item(251, 5).Value = 45380.52083333334
$ws.Cells.Item(251, 6).Value = 'FC Voluntari'
$ws.Cells.Item(251, 7).Value = 'AFC Hermannstadt'
$ws.Cells.Item(251, 8).Value = 1
$ws.Cells.Item(251, 9).Value = 0
$ws.Cells.Item(251, 10).Value = 'H'
$ws.Cells.Item(251, 11).Value = 3
$ws.Cells.Item(251, 12).Value = 3.1
$ws.Cells.Item(251, 13).Value = 2.45
$ws.Cells.Item(251, 14).Value = 3
$ws.Cells.Item(251, 15).Value = 3.1
$ws.Cells.Item(251, 16).Value = 2.45
$ws.Cells.Item(251, 17).Value = 0.25
$ws.Cells.Item(251, 18).Value = 1.775
$ws.Cells.Item(251, 19).Value = 2.1
$ws.Cells.Item(251, 20).Value = 2
$ws.Cells.Item(251, 21).Value = 1.95
$ws.Cells.Item(251, 22).Value = 1.9
$ws.Cells.Item(251, 23).Value = 2
$ws.Cells.Item(251, 24).Value = -1
$ws.Cells.Item(251, 25).Value = -1
$ws.Cells.Item(251, 26).Value = 0.7749999999999999
$ws.Cells.Item(251, 27).Value = -1
$ws.Cells.Item(251, 28).Value = -1
$ws.Cells.Item(251, 29).Value = 0.8999999999999999

# Row 252
$ws.Cells.Item(252, 1).Value = 250
$ws.Cells.Item(252, 2).Value = 7951750
$ws.Cells.Item(252, 3).Value = 'Romania Liga I'
$ws.Cells.Item(252, 4).Value = 'Romania Liga I'
$ws.Cells.Item(252, 5).Value = 45380.64583333334
$ws.Cells.Item(252, 6).Value = 'ACS Sepsi'
$ws.Cells.Item(252, 7).Value = 'CFR Cluj'
$ws.Cells.Item(252, 8).Value = 1
$ws.Cells.Item(252, 9).Value = 1
$ws.Cells.Item(252, 10).Value = 'D'
$ws.Cells.Item(252, 11).Value = 3.3
$ws.Cells.Item(252, 12).Value = 3.4
$ws.Cells.Item(252, 13).Value = 2.15
$ws.Cells.Item(252, 14).Value = 3.4
$ws.Cells.Item(252, 15).Value = 3.5
$ws.Cells.Item(252, 16).Value = 2.05
$ws.Cells.Item(252, 17).Value = 0.25
$ws.Cells.Item(252, 18).Value = 2.05
$ws.Cells.Item(252, 19).Value = 1.8
$ws.Cells.Item(252, 20).Value = 2.5
$ws.Cells.Item(252, 21).Value = 2.05
$ws.Cells.Item(252, 22).Value = 1.8
$ws.Cells.Item(252, 23).Value = -1
$ws.Cells.Item(252, 24).Value = 2.5
$ws.Cells.Item(252, 25).Value = -1
$ws.Cells.Item(252, 26).Value = 0.5249999999999999
$ws.Cells.Item(252, 27).Value = -0.5
$ws.Cells.Item(252, 28).Value = -1
$ws.Cells.Item(252, 29).Value = 0.8

# Row 253
$ws.Cells.Item(253, 1).Value = 251
$ws.Cells.Item(253, 2).Value = 8010912
$ws.Cells.Item(253, 3).Value = 'Romania Liga I'
$ws.Cells.Item(253, 4).Value = 'Romania Liga I'
$ws.Cells.Item(253, 5).Value = 45381.44791666666
$ws.Cells.Item(253, 6).Value = 'FC Botosani'
$ws.Cells.Item(253, 7).Value = 'CSM Politehnica Iasi'
$ws.Cells.Item(253, 8).Value = 2
$ws.Cells.Item(253, 9).Value = 1
$ws.Cells.Item(253, 10).Value = 'H'
$ws.Cells.Item(253, 11).Value = 2.55
$ws.Cells.Item(253, 12).Value = 3.1
$ws.Cells.Item(253, 13).Value = 2.875
$ws.Cells.Item(253, 14).Value = 2.2
$ws.Cells.Item(253, 15).Value = 3.2
$ws.Cells.Item(253, 16).Value = 3.5
$ws.Cells.Item(253, 17).Value = -0.25
$ws.Cells.Item(253, 18).Value = 1.925
$ws.Cells.Item(253, 19).Value = 1.925
$ws.Cells.Item(253, 20).Value = 2
$ws.Cells.Item(253, 21).Value = 1.85
$ws.Cells.Item(253, 22).Value = 2
$ws.Cells.Item(253, 23).Value = 1.2
$ws.Cells.Item(253, 24).Value = -1
$ws.Cells.Item(253, 25).Value = -1
$ws.Cells.Item(253, 26).Value = 0.925
$ws.Cells.Item(253, 27).Value = -1
$ws.Cells.Item(253, 28).Value = 0.8500000000000001
$ws.Cells.Item(253, 29).Value = -1

# Row 254
$ws.Cells.Item(254, 1).Value = 252
$ws.Cells.Item(254, 2).Value = 8010913
$ws.Cells.Item(254, 3).Value = 'Romania Liga I'
$ws.Cells.Item(254, 4).Value = 'Romania Liga I'
$ws.Cells.Item(254, 5).Value = 45381.54166666666
$ws.Cells.Item(254, 6).Value = 'Universitatea Cluj'
$ws.Cells.Item(254, 7).Value = 'ACS UTA Batrana Doamna'
$ws.Cells.Item(254, 8).Value = 0
$ws.Cells.Item(254, 9).Value = 0
$ws.Cells.Item(254, 10).Value = 'D'
$ws.Cells.Item(254, 11).Value = 1.95
$ws.Cells.Item(254, 12).Value = 3.4
$ws.Cells.Item(254, 13).Value = 4
$ws.Cells.Item(254, 14).Value = 1.75
$ws.Cells.Item(254, 15).Value = 3.5
$ws.Cells.Item(254, 16).Value = 4.5
$ws.Cells.Item(254, 17).Value = -0.75
$ws.Cells.Item(254, 18).Value = 2.05
$ws.Cells.Item(254, 19).Value = 1.8
$ws.Cells.Item(254, 20).Value = 2
$ws.Cells.Item(254, 21).Value = 1.75
$ws.Cells.Item(254, 22).Value = 2.05
$ws.Cells.Item(254, 23).Value = -1
$ws.Cells.Item(254, 24).Value = 2.5
$ws.Cells.Item(254, 25).Value = -1
$ws.Cells.Item(254, 26).Value = -1
$ws.Cells.Item(254, 27).Value = 0.8
$ws.Cells.Item(254, 28).Value = -1
$ws.Cells.Item(254, 29).Value = 1.05

# Row 255
$ws.Cells.Item(255, 1).Value = 253
$ws.Cells.Item(255, 2).Value = 7951749
$ws.Cells.Item(255, 3).Value = 'Romania Liga I'
$ws.Cells.Item(255, 4).Value = 'Romania Liga I'
$ws.Cells.Item(255, 5).Value = 45381.66666666666
$ws.Cells.Item(255, 6).Value = 'CS U Craiova'
$ws.Cells.Item(255, 7).Value = 'Rapid Bucuresti'
$ws.Cells.Item(255, 8).Value = 2
$ws.Cells.Item(255, 9).Value = 1
$ws.Cells.Item(255, 10).Value = 'H'
$ws.Cells.Item(255, 11).Value = 2.1
$ws.Cells.Item(255, 12).Value = 3.4
$ws.Cells.Item(255, 13).Value = 3.3
$ws.Cells.Item(255, 14).Value = 2.45
$ws.Cells.Item(255, 15).Value = 3.4
$ws.Cells.Item(255, 16).Value = 2.7
$ws.Cells.Item(255, 17).Value = 0
$ws.Cells.Item(255, 18).Value = 1.875
$ws.Cells.Item(255, 19).Value = 1.975
$ws.Cells.Item(255, 20).Value = 2.5
$ws.Cells.Item(255, 21).Value = 1.875
$ws.Cells.Item(255, 22).Value = 1.975
$ws.Cells.Item(255, 23).Value = 1.45
$ws.Cells.Item(255, 24).Value = -1
$ws.Cells.Item(255, 25).Value = -1
$ws.Cells.Item(255, 26).Value = 0.875
$ws.Cells.Item(255, 27).Value = -1
$ws.Cells.Item(255, 28).Value = 0.875
$ws.Cells.Item(255, 29).Value = -1

# Row 256
$ws.Cells.Item(256, 1).Value = 254
$ws.Cells.Item(256, 2).Value = 7951780
$ws.Cells.Item(256, 3).Value = 'Romania Liga I'
$ws.Cells.Item(256, 4).Value = 'Romania Liga I'
$ws.Cells.Item(256, 5).Value = 45383.60416666666
$ws.Cells.Item(256, 6).Value = 'Dinamo Bucharest'
$ws.Cells.Item(256, 7).Value = 'Petrolul Ploiesti'
$ws.Cells.Item(256, 11).Value = 2.3
$ws.Cells.Item(256, 12).Value = 3
$ws.Cells.Item(256, 13).Value = 3.4
$ws.Cells.Item(256, 14).Value = 2.5
$ws.Cells.Item(256, 15).Value = 3
$ws.Cells.Item(256, 16).Value = 3
$ws.Cells.Item(256, 17).Value = 0
$ws.Cells.Item(256, 18).Value = 1.75
$ws.Cells.Item(256, 19).Value = 2.125
$ws.Cells.Item(256, 20).Value = 2
$ws.Cells.Item(256, 21).Value = 2
$ws.Cells.Item(256, 22).Value = 1.85
$ws.Cells.Item(256, 23).Value = 0
$ws.Cells.Item(256, 24).Value = 0
$ws.Cells.Item(256, 25).Value = 0
$ws.Cells.Item(256, 26).Value = 0
$ws.Cells.Item(256, 27).Value = 0
